# Add four new fynd (observation) rows (13-16) to the Artfynd sheet,
# matching the source OOXML diff. Values that Excel would otherwise
# auto-convert (full ISO dates, and the literal text "1") are forced
# to stay as text via NumberFormat = "@" before assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("A13").Value = 131108314
$ws.Range("B13").Value = 91808
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 1202
$ws.Range("F13").Value = "Ullticka"
$ws.Range("G13").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H13").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P13").Value = "Paljack, Mpd"
$ws.Range("Q13").Value = 601153
$ws.Range("R13").Value = 6977380
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = "Västernorrland"
$ws.Range("U13").Value = "Sundsvall"
$ws.Range("V13").Value = "Medelpad"
$ws.Range("W13").Value = "Liden"
$ws.Range("X13").Value = "2025-1088"
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = "2025-08-27"
$ws.Range("Z13").Value = "08:29"
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = "2025-08-27"
$ws.Range("AB13").Value = "08:29"
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
$ws.Range("AW13").Value = "David Isaksson"
$ws.Range("AX13").Value = "Anders Forsberg"
$ws.Range("AY13").Value = "Kustpaketet"

# Row 14
$ws.Range("A14").Value = 131108277
$ws.Range("B14").Value = 57881
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 100049
$ws.Range("F14").Value = "Spillkråka"
$ws.Range("G14").Value = "Dryocopus martius"
$ws.Range("H14").Value = "(Linnaeus, 1758)"
$ws.Range("P14").Value = "Paljack, Mpd"
$ws.Range("Q14").Value = 601141
$ws.Range("R14").Value = 6977358
$ws.Range("S14").Value = 10
$ws.Range("T14").Value = "Västernorrland"
$ws.Range("U14").Value = "Sundsvall"
$ws.Range("V14").Value = "Medelpad"
$ws.Range("W14").Value = "Liden"
$ws.Range("X14").Value = "2025-1083"
$ws.Range("Y14").NumberFormat = "@"
$ws.Range("Y14").Value = "2025-08-27"
$ws.Range("Z14").Value = "07:45"
$ws.Range("AA14").NumberFormat = "@"
$ws.Range("AA14").Value = "2025-08-27"
$ws.Range("AB14").Value = "07:45"
$ws.Range("AD14").Value = $false
$ws.Range("AE14").Value = $false
$ws.Range("AG14").Value = $false
$ws.Range("AW14").Value = "David Isaksson"
$ws.Range("AX14").Value = "Anders Forsberg"
$ws.Range("AY14").Value = "Kustpaketet"

# Row 15
$ws.Range("A15").Value = 131108166
$ws.Range("B15").Value = 80348
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 6458
$ws.Range("F15").Value = "Lunglav"
$ws.Range("G15").Value = "Lobaria pulmonaria"
$ws.Range("H15").Value = "(L.) Hoffm."
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "1"
$ws.Range("J15").Value = "m²"
$ws.Range("P15").Value = "Paljack, Mpd"
$ws.Range("Q15").Value = 601219
$ws.Range("R15").Value = 6977464
$ws.Range("S15").Value = 10
$ws.Range("T15").Value = "Västernorrland"
$ws.Range("U15").Value = "Sundsvall"
$ws.Range("V15").Value = "Medelpad"
$ws.Range("W15").Value = "Liden"
$ws.Range("X15").Value = "2025-1054"
$ws.Range("Y15").NumberFormat = "@"
$ws.Range("Y15").Value = "2025-08-27"
$ws.Range("Z15").Value = "09:14"
$ws.Range("AA15").NumberFormat = "@"
$ws.Range("AA15").Value = "2025-08-27"
$ws.Range("AB15").Value = "09:14"
$ws.Range("AD15").Value = $false
$ws.Range("AE15").Value = $false
$ws.Range("AG15").Value = $false
$ws.Range("AW15").Value = "David Isaksson"
$ws.Range("AX15").Value = "Elsa Fogelström, Anders Forsberg"
$ws.Range("AY15").Value = "Kustpaketet"

# Row 16
$ws.Range("A16").Value = 131108169
$ws.Range("B16").Value = 91808
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 1202
$ws.Range("F16").Value = "Ullticka"
$ws.Range("G16").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H16").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P16").Value = "Paljack, Mpd"
$ws.Range("Q16").Value = 601162
$ws.Range("R16").Value = 6977363
$ws.Range("S16").Value = 10
$ws.Range("T16").Value = "Västernorrland"
$ws.Range("U16").Value = "Sundsvall"
$ws.Range("V16").Value = "Medelpad"
$ws.Range("W16").Value = "Liden"
$ws.Range("X16").Value = "2025-1089"
$ws.Range("Y16").NumberFormat = "@"
$ws.Range("Y16").Value = "2025-08-27"
$ws.Range("Z16").Value = "08:37"
$ws.Range("AA16").NumberFormat = "@"
$ws.Range("AA16").Value = "2025-08-27"
$ws.Range("AB16").Value = "08:37"
$ws.Range("AD16").Value = $false
$ws.Range("AE16").Value = $false
$ws.Range("AG16").Value = $false
$ws.Range("AW16").Value = "David Isaksson"
$ws.Range("AX16").Value = "Anders Forsberg"
$ws.Range("AY16").Value = "Kustpaketet"


Write-Output "Added rows 13-16 to sheet '$($ws.Name)'"
